$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 66: Cave of Memories (Purgatory)
$ws.Range("A66").Value = "Cave of Memories"
$ws.Range("B66").Value = "Purgatory"
$ws.Range("D66").Value = "Simple Brass Key"
$ws.Range("E66").Value = "A dark, damp cave full of creatures locked away. Some say they are memories that people burry deep within them selves. Memories that twist and become monsterous"
$ws.Range("G66").Value = 1
$ws.Range("J66").Value = 1312
$ws.Range("K66").Value = 560
$ws.Range("L66").Value = 11
$ws.Range("M66").Value = "No"

# Row 67: Alchemcially corrupted graveyard (Delusional Memories)
$ws.Range("A67").Value = "Alchemcially corrupted graveyard"
$ws.Range("B67").Value = "Delusional Memories"
$ws.Range("E67").Value = "A grave yard by the old church. The souls of the children who were corrupted by the churches twisted alchemical practices"
$ws.Range("G67").Value = 1
$ws.Range("I67").Value = 1
$ws.Range("J67").Value = 1280
$ws.Range("K67").Value = 2064
$ws.Range("M67").Value = "No"

# Column A widened to fit the new, longer location name (bestFit recalculation)
$ws.Columns.Item(1).ColumnWidth = 38
